$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999998593648143
$ws.Range("A2").Value = 0.99531893239156999
$ws.Range("A3").Value = 0.98144903717425758
$ws.Range("A4").Value = 0.97521877446479577
$ws.Range("A5").Value = 0.96601031874524734
$ws.Range("A6").Value = 0.94421381987814001
$ws.Range("A7").Value = 0.94469903265270216
$ws.Range("A8").Value = 0.94362476433294651
$ws.Range("A9").Value = 0.94702190820974164
$ws.Range("A10").Value = 0.95140410476153014
$ws.Range("A11").Value = 0.95201357679876042
$ws.Range("A12").Value = 0.95212125043906104
$ws.Range("A13").Value = 0.95423107512129035
$ws.Range("A14").Value = 0.95588857584347187
$ws.Range("A15").Value = 0.95885858986718031
$ws.Range("A16").Value = 0.96300448854724097
$ws.Range("A17").Value = 0.97326398684909488
$ws.Range("A18").Value = 0.97883816579941518
$ws.Range("A19").Value = 0.99312447520142677
$ws.Range("A20").Value = 0.98600787155506198
$ws.Range("A21").Value = 0.98460946664570703
$ws.Range("A22").Value = 0.98334498472844167
$ws.Range("A23").Value = 0.97726710178208909
$ws.Range("A24").Value = 0.97061927872649512
$ws.Range("A25").Value = 0.96416240832850564
$ws.Range("A26").Value = 0.94055899724449232
$ws.Range("A27").Value = 0.93571312125714057
$ws.Range("A28").Value = 0.91424254279953709
$ws.Range("A29").Value = 0.89897242190751991
$ws.Range("A30").Value = 0.89240244695079318
$ws.Range("A31").Value = 0.88474937441165413
$ws.Range("A32").Value = 0.88307018771941825
$ws.Range("A33").Value = 0.88255021291278068
